$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the obsolete row 53 (Caso 6337 / PARAGUAY 4383); this shifts
#    rows 54-62 up to become rows 53-61.
$ws.Rows.Item(53).Delete()

# 2. After the shift, the row that is now 55 (old row 56, Ohiggins 1611)
#    needs its 'Caso' value corrected from -504 to 6506.
$ws.Range("A55").Value2 = "'6506"

# 3. Append the 7 brand-new rows (62-68) reported in this update.

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'6502"
$arr[0,1] = "'7/25/2025"
$arr[0,2] = "CIUDAD DE LA PAZ 1511"
$arr[0,3] = "'13"
$arr[0,4] = "'808571972"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Picada"
$arr[0,8] = 1
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Pasante"
$arr[0,12] = -58.452907
$arr[0,13] = -34.567508
$arr[0,14] = "Colegiales"
$arr[0,15] = "Capital Norte"
$ws.Range("A62:P62").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'6504"
$arr[0,1] = "'7/25/2025"
$arr[0,2] = "CIUDAD DE LA PAZ 1278"
$arr[0,3] = "'13"
$arr[0,4] = "'808571974"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Inclinada"
$arr[0,8] = 1
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Pasante"
$arr[0,12] = -58.450753
$arr[0,13] = -34.5688
$arr[0,14] = "Colegiales"
$arr[0,15] = "Capital Norte"
$ws.Range("A63:P63").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'6512"
$arr[0,1] = "'7/28/2025"
$arr[0,2] = "GASCON 1195"
$arr[0,3] = "'14"
$arr[0,4] = "'808571975"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Picada"
$arr[0,8] = 1
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Pasante"
$arr[0,12] = -58.423127
$arr[0,13] = -34.596476
$arr[0,14] = "Palermo"
$arr[0,15] = "Capital Sur"
$ws.Range("A64:P64").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'6513"
$arr[0,1] = "'7/28/2025"
$arr[0,2] = "DORREGO 1925"
$arr[0,3] = "'14"
$arr[0,4] = "'808571976"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Picada"
$arr[0,8] = 1
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Pasante"
$arr[0,12] = -58.441281
$arr[0,13] = -34.579867
$arr[0,14] = "Palermo"
$arr[0,15] = "Capital Sur"
$ws.Range("A65:P65").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'6519"
$arr[0,1] = "'7/28/2025"
$arr[0,2] = "SALGUERO, JERONIMO 2874"
$arr[0,3] = "'14"
$arr[0,4] = "'808571977"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Reparar rienda"
$arr[0,8] = 1
$arr[0,9] = "Tensor"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Terminal"
$arr[0,12] = -58.407256
$arr[0,13] = -34.578976
$arr[0,14] = "Palermo"
$arr[0,15] = "Capital Sur"
$ws.Range("A66:P66").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'-534"
$arr[0,1] = "'7/28/2025"
$arr[0,2] = "Jose Aaron Salmun Feijoo 325"
$arr[0,3] = "'4"
$arr[0,4] = "'808571999"
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Cambiar terminal"
$arr[0,8] = 0
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Terminal"
$arr[0,12] = -58.3793
$arr[0,13] = -34.636079
$arr[0,14] = "San Telmo"
$arr[0,15] = "Capital Sur"
$ws.Range("A67:P67").Value2 = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "'-535"
$arr[0,1] = "'7/28/2025"
$arr[0,2] = "Jose Aaron Salmun Feijoo 363"
$arr[0,3] = "'4"
$arr[0,4] = $null
$arr[0,5] = "NEW"
$arr[0,6] = "Pendiente"
$arr[0,7] = "Colocar prfv pasante"
$arr[0,8] = 0
$arr[0,9] = "Cambio"
$arr[0,10] = "Sin equipos"
$arr[0,11] = "Pasante"
$arr[0,12] = -58.379294
$arr[0,13] = -34.636313
$arr[0,14] = "San Telmo"
$arr[0,15] = "Capital Sur"
$ws.Range("A68:P68").Value2 = $arr
